# Validators.xlsx — "Started to work on the type-extensions" edit.
#
# The real content change here is on the "Functions" sheet: the example
# value used to illustrate the String/Number/Boolean/Array columns for the
# IsNotNull/IsNull (row 2) and IsUndefined/IsDefined (row 3) rows is updated
# from the placeholder "X" to "Y- Test" (and, for F3, "Y-Test"). That also
# moves the window's selection from J20 to K15 and scrolls the sheet so
# column F is back at the left edge.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")
$ws.Activate() | Out-Null

# Row 2 (IsNotNull / IsNull): String, Number, Boolean and Array example
# cells all become "Y- Test".
$ws.Range("C2:F2").Value = "Y- Test"

# Row 3 (IsUndefined / IsDefined): String/Number/Boolean become "Y- Test",
# but the Array example (F3) becomes "Y-Test" (no space before "Test").
$ws.Range("C3:E3").Value = "Y- Test"
$ws.Range("F3").Value = "Y-Test"

# Scroll the view so column F sits at the left edge (topLeftCell F1) and
# move the active selection to K15.
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K15").Select() | Out-Null
